$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the "To" and "CC" user values in row 2 with the automation test users
$ws.Range("A2").Value = "AutoTestAdmin"
$ws.Range("B2").Value = "AutoTestUser"
